$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New matchup rows to append (Player_1, Points_1, Player_2, Points_2)
$rows = @(
    @(5,1,4,2),
    @(5,0,7,3),
    @(5,0,5,2),
    @(6,0,6,2),
    @(4,0,4,3),
    @(4,1,4,2),
    @(4,1,4,2),
    @(5,2,3,1),
    @(3,2,4,1),
    @(6,0,4,3),
    @(3,3,3,0),
    @(3,2,4,0),
    @(7,0,5,2),
    @(4,1,4,2),
    @(3,2,3,1),
    @(6,2,6,1),
    @(6,2,6,0),
    @(2,2,4,1),
    @(6,2,7,1),
    @(4,2,3,1),
    @(6,2,5,1),
    @(5,0,5,2),
    @(2,1,2,2)
)

$startRow = 1525
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
}

$lastRow = $startRow + $rows.Count - 1
$nextRow = $lastRow + 1

$ws.Range("A" + $nextRow).Select()
